# Insert a new daily data row right before the existing row 880
# (2026/12/29 火 ...), pushing all subsequent rows down by one.
# The new row carries the next reading for 2026/02/25 (水), continuing
# the C-column sequence already present in rows 878-879 (1, 5, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 880:921 down to 881:922, creating a blank row 880.
$ws.Rows.Item(880).Insert()

# Column A holds a date formatted as plain text (e.g. "2026/02/25"), not
# a real Excel date. Force text formatting before assigning the value so
# Excel doesn't auto-convert the string into a date serial number, then
# drop back to the default (unstyled) cell style to match the rest of
# the sheet's data rows.
$ws.Range("A880").NumberFormat = "@"
$ws.Range("A880").Value = "2026/02/25"
$ws.Range("A880").Style = "Normal"

$ws.Range("B880").Value = "水"
$ws.Range("C880").Value = 8
$ws.Range("D880").Value = 201
